$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 30; this shifts existing rows 30-141 down to 31-142
$ws.Rows.Item(30).Insert()

# Populate the newly inserted row 30 with the new record's data
$ws.Cells.Item(30, 1).Value = 5
$ws.Cells.Item(30, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(30, 3).Value = "Maule"
$ws.Cells.Item(30, 4).Value = 44487
$ws.Cells.Item(30, 5).Value = 7
$ws.Cells.Item(30, 6).Value = 100112045
$ws.Cells.Item(30, 7).Value = "Zapallo"
$ws.Cells.Item(30, 8).Value = "Paine"
$ws.Cells.Item(30, 9).Value = "1a (guarda)"
$ws.Cells.Item(30, 10).Value = 3000
$ws.Cells.Item(30, 11).Value = 100
$ws.Cells.Item(30, 12).Value = 100
$ws.Cells.Item(30, 13).Value = 100
$ws.Cells.Item(30, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(30, 15).Value = "Región del Maule"
$ws.Cells.Item(30, 16).Value = 100
$ws.Cells.Item(30, 17).Value = 1
$ws.Cells.Item(30, 18).Value = "Hortaliza"
